$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows before the current row 113 (the old "Dulce o Americano"
# entry gets pushed down to row 115), then populate the two new rows with the
# new "Choclero" price entries dated 44911.
$ws.Rows("113:114").Insert()

$ws.Range("A113").Value = 11
$ws.Range("B113").Value = "Vega Monumental Concepción"
$ws.Range("C113").Value = "Bíobío"
$ws.Range("D113").Value = 44911
$ws.Range("E113").Value = 8
$ws.Range("F113").Value = 100112024
$ws.Range("G113").Value = "Choclo"
$ws.Range("H113").Value = "Choclero"
$ws.Range("I113").Value = "Primera"
$ws.Range("J113").Value = 15000
$ws.Range("K113").Value = 300
$ws.Range("L113").Value = 350
$ws.Range("M113").Value = 333
$ws.Range("N113").Value = "`$/unidad"
$ws.Range("O113").Value = "Región de O'Higgins"
$ws.Range("P113").Value = 333
$ws.Range("Q113").Value = 1
$ws.Range("R113").Value = "Hortaliza"

$ws.Range("A114").Value = 11
$ws.Range("B114").Value = "Vega Monumental Concepción"
$ws.Range("C114").Value = "Bíobío"
$ws.Range("D114").Value = 44911
$ws.Range("E114").Value = 8
$ws.Range("F114").Value = 100112024
$ws.Range("G114").Value = "Choclo"
$ws.Range("H114").Value = "Choclero"
$ws.Range("I114").Value = "Segunda"
$ws.Range("J114").Value = 5000
$ws.Range("K114").Value = 250
$ws.Range("L114").Value = 250
$ws.Range("M114").Value = 250
$ws.Range("N114").Value = "`$/unidad"
$ws.Range("O114").Value = "Región de O'Higgins"
$ws.Range("P114").Value = 250
$ws.Range("Q114").Value = 1
$ws.Range("R114").Value = "Hortaliza"
